$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 75
$ws.Range("A75").Value = 111949575
$ws.Range("B75").Value = 96348
$ws.Range("C75").Value = "Ovaliderad"
$ws.Range("D75").Value = "VU"
$ws.Range("E75").Value = 220787
$ws.Range("F75").Value = "Knärot"
$ws.Range("G75").Value = "Goodyera repens"
$ws.Range("H75").Value = "(L.) R. Br."
$ws.Range("I75").NumberFormat = "@"
$ws.Range("I75").Value = "15"
$ws.Range("I75").ClearFormats()
$ws.Range("J75").Value = "plantor/tuvor"
# K75 left blank (source is an empty inline string cell)
$ws.Range("P75").Value = "Slättesmyran (Slättesmyran), Ång"
$ws.Range("Q75").Value = 580471.3517951096
$ws.Range("R75").Value = 7053333.257918903
$ws.Range("S75").Value = 1
$ws.Range("T75").Value = "Västernorrland"
$ws.Range("U75").Value = "Sollefteå"
$ws.Range("V75").Value = "Ångermanland"
$ws.Range("W75").Value = "Ramsele"
$ws.Range("Y75").NumberFormat = "@"
$ws.Range("Y75").Value = "2023-09-07"
$ws.Range("Y75").ClearFormats()
$ws.Range("Z75").Value = "19:05"
$ws.Range("AA75").NumberFormat = "@"
$ws.Range("AA75").Value = "2023-09-07"
$ws.Range("AA75").ClearFormats()
$ws.Range("AB75").Value = "19:05"
$ws.Range("AD75").Value = $false
$ws.Range("AE75").Value = $false
$ws.Range("AG75").Value = $false
# AT75 left blank (source is an empty inline string cell)
$ws.Range("AW75").Value = "Kamilla Andersson"
$ws.Range("AX75").Value = "Kamilla Andersson"
# AY75 left blank (source is an empty inline string cell)

# Row 76
$ws.Range("A76").Value = 111949678
$ws.Range("B76").Value = 96348
$ws.Range("C76").Value = "Ovaliderad"
$ws.Range("D76").Value = "VU"
$ws.Range("E76").Value = 220787
$ws.Range("F76").Value = "Knärot"
$ws.Range("G76").Value = "Goodyera repens"
$ws.Range("H76").Value = "(L.) R. Br."
$ws.Range("I76").NumberFormat = "@"
$ws.Range("I76").Value = "7"
$ws.Range("I76").ClearFormats()
# K76 left blank (source is an empty inline string cell)
$ws.Range("P76").Value = "Slättesmyran (Slättesmyran), Ång"
$ws.Range("Q76").Value = 580467.4207067642
$ws.Range("R76").Value = 7053330.04139028
$ws.Range("S76").Value = 2
$ws.Range("T76").Value = "Västernorrland"
$ws.Range("U76").Value = "Sollefteå"
$ws.Range("V76").Value = "Ångermanland"
$ws.Range("W76").Value = "Ramsele"
$ws.Range("Y76").NumberFormat = "@"
$ws.Range("Y76").Value = "2023-09-07"
$ws.Range("Y76").ClearFormats()
$ws.Range("Z76").Value = "19:11"
$ws.Range("AA76").NumberFormat = "@"
$ws.Range("AA76").Value = "2023-09-07"
$ws.Range("AA76").ClearFormats()
$ws.Range("AB76").Value = "19:11"
$ws.Range("AD76").Value = $false
$ws.Range("AE76").Value = $false
$ws.Range("AG76").Value = $false
# AT76 left blank (source is an empty inline string cell)
$ws.Range("AW76").Value = "Kim Hultgren"
$ws.Range("AX76").Value = "Kim Hultgren"
# AY76 left blank (source is an empty inline string cell)

# Row 77
$ws.Range("A77").Value = 111949591
$ws.Range("B77").Value = 96348
$ws.Range("C77").Value = "Ovaliderad"
$ws.Range("D77").Value = "VU"
$ws.Range("E77").Value = 220787
$ws.Range("F77").Value = "Knärot"
$ws.Range("G77").Value = "Goodyera repens"
$ws.Range("H77").Value = "(L.) R. Br."
$ws.Range("I77").NumberFormat = "@"
$ws.Range("I77").Value = "20"
$ws.Range("I77").ClearFormats()
# K77 left blank (source is an empty inline string cell)
$ws.Range("P77").Value = "Slättesmyran (Slättesmyran), Ång"
$ws.Range("Q77").Value = 580476.1122211452
$ws.Range("R77").Value = 7053321.356648902
$ws.Range("S77").Value = 2
$ws.Range("T77").Value = "Västernorrland"
$ws.Range("U77").Value = "Sollefteå"
$ws.Range("V77").Value = "Ångermanland"
$ws.Range("W77").Value = "Ramsele"
$ws.Range("Y77").NumberFormat = "@"
$ws.Range("Y77").Value = "2023-09-07"
$ws.Range("Y77").ClearFormats()
$ws.Range("Z77").Value = "19:07"
$ws.Range("AA77").NumberFormat = "@"
$ws.Range("AA77").Value = "2023-09-07"
$ws.Range("AA77").ClearFormats()
$ws.Range("AB77").Value = "19:07"
$ws.Range("AD77").Value = $false
$ws.Range("AE77").Value = $false
$ws.Range("AG77").Value = $false
# AT77 left blank (source is an empty inline string cell)
$ws.Range("AW77").Value = "Kim Hultgren"
$ws.Range("AX77").Value = "Kim Hultgren"
# AY77 left blank (source is an empty inline string cell)

# Row 78
$ws.Range("A78").Value = 111950184
$ws.Range("B78").Value = 56543
$ws.Range("C78").Value = "Ovaliderad"
$ws.Range("D78").Value = "NT"
$ws.Range("E78").Value = 103021
$ws.Range("F78").Value = "Talltita"
$ws.Range("G78").Value = "Poecile montanus"
$ws.Range("H78").Value = "(Conrad von Baldenstein, 1827)"
# I78 left blank (source is an empty inline string cell)
# K78 left blank (source is an empty inline string cell)
$ws.Range("P78").Value = "Slättesmyran (Slättesmyran), Ång"
$ws.Range("Q78").Value = 580446.7330953531
$ws.Range("R78").Value = 7053301.910512885
$ws.Range("S78").Value = 10
$ws.Range("T78").Value = "Västernorrland"
$ws.Range("U78").Value = "Sollefteå"
$ws.Range("V78").Value = "Ångermanland"
$ws.Range("W78").Value = "Ramsele"
$ws.Range("Y78").NumberFormat = "@"
$ws.Range("Y78").Value = "2023-09-07"
$ws.Range("Y78").ClearFormats()
$ws.Range("Z78").Value = "19:37"
$ws.Range("AA78").NumberFormat = "@"
$ws.Range("AA78").Value = "2023-09-07"
$ws.Range("AA78").ClearFormats()
$ws.Range("AB78").Value = "19:37"
$ws.Range("AD78").Value = $false
$ws.Range("AE78").Value = $false
$ws.Range("AG78").Value = $false
# AT78 left blank (source is an empty inline string cell)
$ws.Range("AW78").Value = "Kim Hultgren"
$ws.Range("AX78").Value = "Kim Hultgren"
# AY78 left blank (source is an empty inline string cell)

# Row 79
$ws.Range("A79").Value = 111949317
$ws.Range("B79").Value = 96265
$ws.Range("C79").Value = "Ovaliderad"
$ws.Range("D79").Value = "LC"
$ws.Range("E79").Value = 219790
$ws.Range("F79").Value = "Fläcknycklar"
$ws.Range("G79").Value = "Dactylorhiza maculata"
$ws.Range("H79").Value = "(L.) Soó"
# I79 left blank (source is an empty inline string cell)
# K79 left blank (source is an empty inline string cell)
$ws.Range("P79").Value = "Slättesmyran (Slättesmyran), Ång"
$ws.Range("Q79").Value = 580500.003505226
$ws.Range("R79").Value = 7053328.641698814
$ws.Range("S79").Value = 2
$ws.Range("T79").Value = "Västernorrland"
$ws.Range("U79").Value = "Sollefteå"
$ws.Range("V79").Value = "Ångermanland"
$ws.Range("W79").Value = "Ramsele"
$ws.Range("Y79").NumberFormat = "@"
$ws.Range("Y79").Value = "2023-09-07"
$ws.Range("Y79").ClearFormats()
$ws.Range("Z79").Value = "18:54"
$ws.Range("AA79").NumberFormat = "@"
$ws.Range("AA79").Value = "2023-09-07"
$ws.Range("AA79").ClearFormats()
$ws.Range("AB79").Value = "18:54"
$ws.Range("AD79").Value = $false
$ws.Range("AE79").Value = $false
$ws.Range("AG79").Value = $false
# AT79 left blank (source is an empty inline string cell)
$ws.Range("AW79").Value = "Kim Hultgren"
$ws.Range("AX79").Value = "Kim Hultgren"
# AY79 left blank (source is an empty inline string cell)
